$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Steering Torque candidate" note) now wraps onto more lines,
# so its row height needs to grow to fit the text.
$ws.Rows(11).RowHeight = 36

# Move the view / selection: scroll back to the top of the sheet and
# select E2 instead of the previous C16 / A7 top-left scroll position.
$ws.Range("E2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
